# Audits.xlsx - update the RTM review comment (row 24): it has now been
# closed, so update its close date and state, and record the closer in
# column H (matching the other closed rows). Also move the active
# selection to H23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 24 ("RTM must be reviewed against SRS.") - close date moved from
# 15/5/2022 to 22/5/2022, and state changed from "open" to "closed".
$ws.Range("D24").Value = "22/5/2022"
$ws.Range("F24").Value = "closed"

# Row 24 now also has a value in column H (like the rows above it).
$ws.Range("H24").Value = "25/5/2022"

# Update the active sheet selection to H23.
$ws.Activate() | Out-Null
$ws.Range("H23").Select() | Out-Null
